$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 229 - existing rows 229-261 shift down to 233-265.
$ws.Range("A229:R232").EntireRow.Insert()

# Columns that stay constant across every data row in this block.
$constA = 11
$constB = "Vega Monumental Concepción"
$constC = "Bíobío"
$constE = 8
$constF = 100112004
$constG = "Cebolla"
$constN = "`$/malla 18 kilos"
$constQ = 18
$constR = "Hortaliza"

# New rows of data (matching the diff's target content for rows 229-232).
$rows = @(
    @{ Row = 229; D = 44476; H = "Morada(o)";        I = "1a (cosecha)"; J = 200; K = 7500; L = 8000; M = 7750; O = "Región de Arica y Parinacota"; P = 431 },
    @{ Row = 230; D = 44476; H = "Morada(o)";        I = "2a (cosecha)"; J = 100; K = 7000; L = 7000; M = 7000; O = "Región de Arica y Parinacota"; P = 389 },
    @{ Row = 231; D = 44476; H = "Sin especificar";  I = "1a (guarda)";  J = 600; K = 6000; L = 6500; M = 6250; O = "Región de O'Higgins";          P = 347 },
    @{ Row = 232; D = 44384; H = "Sin especificar";  I = "2a (guarda)";  J = 300; K = 5500; L = 5500; M = 5500; O = "Región de O'Higgins";          P = 306 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $constA
    $ws.Cells.Item($row, 2).Value = $constB
    $ws.Cells.Item($row, 3).Value = $constC
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $constE
    $ws.Cells.Item($row, 6).Value = $constF
    $ws.Cells.Item($row, 7).Value = $constG
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $constN
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $constQ
    $ws.Cells.Item($row, 18).Value = $constR
}
